$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new paragraph before "Exit Visual Studio in the Cloud_SDK
#    solution." containing the new instruction text, with a collapsed
#    "_GoBack" bookmark right after the text (this also implicitly removes
#    the document's pre-existing _GoBack bookmark, since bookmark names must
#    be unique - matching the diff, which drops the old one further down).
# ---------------------------------------------------------------------------
$exitPara = $d.Paragraphs(4)
$exitPara.Range.InsertParagraphBefore()

$newPara = $d.Paragraphs(4)
$newRange = $newPara.Range
# Append a throw-away trailing marker character so the bookmark insertion
# point below is never the very last character of the paragraph (a
# collapsed range placed there resolves incorrectly), then strip it back out.
$newRange.InsertBefore("Inspect the Sample-Live-Sync directory hierarchy and remove and *.pdb and *.zip files, and any other non-public files.X")

$afterText = $d.Paragraphs(4).Range
$bookmarkPos = $afterText.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$afterBookmark = $d.Paragraphs(4).Range
$markerPos = $afterBookmark.End - 2
$markerRange = $d.Range($markerPos, $markerPos + 1)
$markerRange.Delete()

# ---------------------------------------------------------------------------
# 2) Replace "20130620" + the old _GoBack bookmark + "Release0.1.5" with the
#    literal text "<BranchName>.", while keeping the preceding "-ff   " run
#    intact as its own run.
# ---------------------------------------------------------------------------
$branchRange = $d.Content
$branchRange.Find.Execute("20130620Release0.1.5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$branchRange.Delete()
$branchRange.InsertAfter("<BranchName>.")
